# Apply the "Add files via upload" edit: insert a new "实验语料" column,
# fill in the remaining cells of row 4, and append four new bibliography
# rows (5-8) with their associated metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before the old column E ("缺点"), which
#        becomes column F. The new column E is "实验语料". Inserting this
#        way carries the wrap-text style already present on E2/E3 along
#        for the ride, exactly like the target file.
$ws.Columns("E:E").Insert()

$ws.Range("E1").Value = "实验语料"

# --- 2. Finish populating row 4 (D4/E4/F4 were blank before the edit).
$ws.Range("D4").Value = "提出了一种综合了规则、字典与统计学习的方法。`n1.StandfordNER和FudanNER识别出公司名`n2.扩展字典。利用规则生成一个简称字典`n3.对识别出的简称，利用规则生成其对应的简称，判断上下文中是否存在简称，如果存在，则认为识别出了一个全称-简称对，"
$ws.Range("E4").Value = "用了150篇金融新闻报道，通过手工标注出机构名和简称。`n"
$ws.Range("F4").Value = "其实根本的问题有两个：1.全称的识别依赖以StandfordNER和FudanNER`n2.简称的生成依赖于规则"

$ws.Range("D4").WrapText = $true
$ws.Range("E4").WrapText = $true
$ws.Range("F4").WrapText = $true

# --- 3. Row 5 - "中文组织机构名称与简称的识别"
$ws.Range("A5").Value = "中文组织机构名称与简称的识别"
$ws.Range("B5").Value = 2007
$ws.Range("C5").Value = "期刊"
$ws.Range("D5").Value = "全称与简称的识别都基于规则。首先识别出全称，再根据规则生成简称，然后在上下文中检索简称，如果存在，则认为存在一个匹配"
$ws.Range("E5").Value = "含有654个机构名称的280篇文章作为开放测试集"
$ws.Range("D5").WrapText = $true

# --- 4. Row 6 - CARD paper
$ws.Range("A6").Value = "A long journey to short abbreviations: developing an open-source framework for clinical abbreviation recognition and disambiguation (CARD)"
$ws.Range("B6").Value = 2017
$ws.Range("C6").Value = "期刊"
$ws.Range("D6").Value = "提出了一种医疗实体简称（疾病名、药名）识别与消歧的框架(CARD)。对简称的识别视为一个二分类的为题，采用随机森林、决策树、SVM。"
$ws.Range("A6").WrapText = $true
$ws.Range("D6").WrapText = $true

# --- 5. Row 7 - Collaborative Recognition and Recovery ...
$ws.Range("A7").Value = "Collaborative Recognition and `nRecovery of the Chinese Intercept Abbreviation"
$ws.Range("B7").Value = 2017
$ws.Range("D7").Value = "定义的特征还是可以参考一下的`n关于简称的三个方向：recognition、prediction、recovery`n关于CRF公式的引用`n论文里的整个方法还是可以参考一下的，和我的思路挺接近的"
$ws.Range("A7").WrapText = $true
$ws.Range("D7").WrapText = $true

# --- 6. Row 8 - trailing blank row (kept the wrap-text style only).
$ws.Range("A8").WrapText = $true

# --- 7. Row heights to match the re-flowed content.
$ws.Rows("2:2").RowHeight = 57
$ws.Rows("3:3").RowHeight = 28.5
$ws.Rows("4:4").RowHeight = 85.5
$ws.Rows("5:5").RowHeight = 42.75
$ws.Rows("6:6").RowHeight = 71.25
$ws.Rows("7:7").RowHeight = 85.5

# --- 8. Column widths (closest achievable values given the engine's
#        7-pixel Maximum-Digit-Width rounding of the ColumnWidth property).
$ws.Columns("A:A").ColumnWidth = 29.428571428571427
$ws.Columns("B:B").ColumnWidth = 29
$ws.Columns("C:C").ColumnWidth = 20.857142857142858
$ws.Columns("D:D").ColumnWidth = 42.714285714285715
$ws.Columns("E:E").ColumnWidth = 42.714285714285715
$ws.Columns("F:F").ColumnWidth = 58

# --- 9. Selection moves to E4 after the edit.
$null = $ws.Range("E4").Select()
